$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -12.44060000000001
$ws.Range("B10").Value = 5.5795
$ws.Range("B12").Value = 5.517799999999998
$ws.Range("C12").Value = -13.6759
$ws.Range("D12").Value = -7.864000000000001
$ws.Range("D13").Value = -8.643799999999993
$ws.Range("C17").Value = -13.31129999999999
$ws.Range("B18").Value = 6.552599999999996
$ws.Range("D21").Value = -8.527699999999996
$ws.Range("C26").Value = -12.73520000000001
$ws.Range("C27").Value = -13.0246
$ws.Range("C28").Value = -14.0015
$ws.Range("D36").Value = -7.561000000000002
$ws.Range("B37").Value = 9.015699999999995
$ws.Range("C37").Value = -13.05469999999999
$ws.Range("D38").Value = -8.219900000000001
$ws.Range("D41").Value = -8.213599999999992
$ws.Range("D52").Value = -7.776899999999999
$ws.Range("B55").Value = 6.184999999999994
$ws.Range("C65").Value = -12.681
$ws.Range("D67").Value = -7.100199999999997
$ws.Range("B68").Value = 4.710099999999994
$ws.Range("C73").Value = -11.262
$ws.Range("B77").Value = 9.218400000000008
$ws.Range("B78").Value = 10.0124
$ws.Range("C84").Value = -13.4671
$ws.Range("C85").Value = -13.2711
$ws.Range("D89").Value = -8.305099999999999
$ws.Range("C93").Value = -10.315
$ws.Range("C95").Value = -12.90269999999999
$ws.Range("D95").Value = -7.752900000000003
$ws.Range("C98").Value = -12.86910000000001
$ws.Range("C99").Value = -11.25950000000001
$ws.Range("C101").Value = -13.05570000000001
$ws.Range("D105").Value = -7.910700000000003
